$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Update the query text in B2 (ParticipantsTab query)
$newQuery = "MATCH (p:participant)-->(s:study)`nOPTIONAL MATCH (samp:sample)-->(p)`nOPTIONAL MATCH (p)<--(diag:diagnosis)`nOPTIONAL MATCH (samp)<--(f:file)`nOPTIONAL MATCH (f)<--(g:genomic_info)`nWITH s, p, samp, f, g, diag`nWHERE g.library_source in ['Single Cell']`nWITH p`nOPTIONAL MATCH (p)-->(s:study)`nOPTIONAL MATCH (samp:sample)-->(p)`nWITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp`nRETURN `ncoalesce(p.participant_id,'') as ``Participant ID```,`ncoalesce(s.study_name, '') as ``Study Name```,`ncoalesce(s.phs_accession,'') as ``Accession```,`ncoalesce(p.gender,'') as ``Gender```,`ncoalesce(apoc.text.join(samp, ','), '') as ``Samples```nORDER BY p.participant_id`nLIMIT 100"

$ws.Range("B2").Value = $newQuery

# Increase font size for the whole sheet's default style (columns) to 15
$ws.Cells.Font.Size = 15

# Set column-level formatting (style) for columns A:E and beyond, font size 15, no wrap by default
$ws.Columns("A:E").Font.Size = 15
$ws.Columns("A:E").WrapText = $false

# B and C columns keep wrap text (used for long query text)
$ws.Columns("B:C").WrapText = $true

# Row heights - set to autofit appropriate sizes matching target
$ws.Rows(1).RowHeight = 19.5
$ws.Rows(2).RowHeight = 390
$ws.Rows(3).RowHeight = 292.5
$ws.Rows(4).RowHeight = 292.5
$ws.Rows(5).RowHeight = 19.5
$ws.Rows(6).RowHeight = 19.5

$ws.Range("B5").WrapText = $true
$ws.Range("C5").WrapText = $true
$ws.Range("C6").WrapText = $true

# Update selection to E2
$ws.Range("E2").Select()

# Update workbook window position
$excel.Left = -120
$excel.Top = -120
